$wb = $excel.ActiveWorkbook
$count = $wb.Worksheets.Count
$last = $wb.Worksheets.Item($count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$ws.Name = "trend"

# Header row (row 1): shared-string header labels
$headers = @(
  'ANCOM-BC2 (No Filter).30, 0.05 (N = 100)',
  'ANCOM-BC2 (SS Filter).30, 0.05 (N = 100)',
  'ANCOM-BC2 (No Filter).30, 0.1 (N = 100)',
  'ANCOM-BC2 (SS Filter).30, 0.1 (N = 100)',
  'ANCOM-BC2 (No Filter).30, 0.2 (N = 100)',
  'ANCOM-BC2 (SS Filter).30, 0.2 (N = 100)',
  'ANCOM-BC2 (No Filter).30, 0.5 (N = 100)',
  'ANCOM-BC2 (SS Filter).30, 0.5 (N = 100)',
  'ANCOM-BC2 (No Filter).30, 0.9 (N = 100)',
  'ANCOM-BC2 (SS Filter).30, 0.9 (N = 100)',
  'ANCOM-BC2 (No Filter).60, 0.05 (N = 100)',
  'ANCOM-BC2 (SS Filter).60, 0.05 (N = 100)',
  'ANCOM-BC2 (No Filter).60, 0.1 (N = 100)',
  'ANCOM-BC2 (SS Filter).60, 0.1 (N = 100)',
  'ANCOM-BC2 (No Filter).60, 0.2 (N = 100)',
  'ANCOM-BC2 (SS Filter).60, 0.2 (N = 100)',
  'ANCOM-BC2 (No Filter).60, 0.5 (N = 100)',
  'ANCOM-BC2 (SS Filter).60, 0.5 (N = 100)',
  'ANCOM-BC2 (No Filter).60, 0.9 (N = 100)',
  'ANCOM-BC2 (SS Filter).60, 0.9 (N = 100)',
  'ANCOM-BC2 (No Filter).90, 0.05 (N = 100)',
  'ANCOM-BC2 (SS Filter).90, 0.05 (N = 100)',
  'ANCOM-BC2 (No Filter).90, 0.1 (N = 100)',
  'ANCOM-BC2 (SS Filter).90, 0.1 (N = 100)',
  'ANCOM-BC2 (No Filter).90, 0.2 (N = 100)',
  'ANCOM-BC2 (SS Filter).90, 0.2 (N = 100)',
  'ANCOM-BC2 (No Filter).90, 0.5 (N = 100)',
  'ANCOM-BC2 (SS Filter).90, 0.5 (N = 100)',
  'ANCOM-BC2 (No Filter).90, 0.9 (N = 100)',
  'ANCOM-BC2 (SS Filter).90, 0.9 (N = 100)',
  'ANCOM-BC2 (No Filter).150, 0.05 (N = 100)',
  'ANCOM-BC2 (SS Filter).150, 0.05 (N = 100)',
  'ANCOM-BC2 (No Filter).150, 0.1 (N = 100)',
  'ANCOM-BC2 (SS Filter).150, 0.1 (N = 100)',
  'ANCOM-BC2 (No Filter).150, 0.2 (N = 100)',
  'ANCOM-BC2 (SS Filter).150, 0.2 (N = 100)',
  'ANCOM-BC2 (No Filter).150, 0.5 (N = 100)',
  'ANCOM-BC2 (SS Filter).150, 0.5 (N = 100)',
  'ANCOM-BC2 (No Filter).150, 0.9 (N = 100)',
  'ANCOM-BC2 (SS Filter).150, 0.9 (N = 100)',
  'ANCOM-BC2 (No Filter).300, 0.05 (N = 100)',
  'ANCOM-BC2 (SS Filter).300, 0.05 (N = 100)',
  'ANCOM-BC2 (No Filter).300, 0.1 (N = 100)',
  'ANCOM-BC2 (SS Filter).300, 0.1 (N = 100)',
  'ANCOM-BC2 (No Filter).300, 0.2 (N = 100)',
  'ANCOM-BC2 (SS Filter).300, 0.2 (N = 100)',
  'ANCOM-BC2 (No Filter).300, 0.5 (N = 100)',
  'ANCOM-BC2 (SS Filter).300, 0.5 (N = 100)',
  'ANCOM-BC2 (No Filter).300, 0.9 (N = 100)',
  'ANCOM-BC2 (SS Filter).300, 0.9 (N = 100)'
)
$headerRange = $ws.Range("A1:AX1")
$headerArr = New-Object 'object[,]' 1,50
for ($i = 0; $i -lt 50; $i++) { $headerArr[0,$i] = $headers[$i] }
$headerRange.Value = $headerArr

# Data rows 2-5: 50 numeric columns each
$row2 = New-Object 'object[,]' 1,50
$row2[0,0] = 1
$row2[0,1] = 0.99
$row2[0,2] = 1
$row2[0,3] = 0.99
$row2[0,4] = 1
$row2[0,5] = 0.99
$row2[0,6] = 0.98
$row2[0,7] = 0.95
$row2[0,8] = 0.46
$row2[0,9] = 0.38
$row2[0,10] = 1
$row2[0,11] = 1
$row2[0,12] = 1
$row2[0,13] = 1
$row2[0,14] = 1
$row2[0,15] = 1
$row2[0,16] = 0.94
$row2[0,17] = 0.9
$row2[0,18] = 0.51
$row2[0,19] = 0.43
$row2[0,20] = 1
$row2[0,21] = 1
$row2[0,22] = 1
$row2[0,23] = 1
$row2[0,24] = 1
$row2[0,25] = 1
$row2[0,26] = 0.95
$row2[0,27] = 0.9
$row2[0,28] = 0.52
$row2[0,29] = 0.44
$row2[0,30] = 1
$row2[0,31] = 1
$row2[0,32] = 1
$row2[0,33] = 1
$row2[0,34] = 1
$row2[0,35] = 1
$row2[0,36] = 0.97
$row2[0,37] = 0.93
$row2[0,38] = 0.53
$row2[0,39] = 0.46
$row2[0,40] = 1
$row2[0,41] = 1
$row2[0,42] = 1
$row2[0,43] = 1
$row2[0,44] = 1
$row2[0,45] = 1
$row2[0,46] = 0.99
$row2[0,47] = 0.98
$row2[0,48] = 0.53
$row2[0,49] = 0.46
$ws.Range("A2:AX2").Value = $row2

$row3 = New-Object 'object[,]' 1,50
$row3[0,0] = 0
$row3[0,1] = 0.03
$row3[0,2] = 0
$row3[0,3] = 0.02
$row3[0,4] = 0
$row3[0,5] = 0.02
$row3[0,6] = 0.04
$row3[0,7] = 0.07
$row3[0,8] = 0.04
$row3[0,9] = 0.06
$row3[0,10] = 0
$row3[0,11] = 0.01
$row3[0,12] = 0
$row3[0,13] = 0.01
$row3[0,14] = 0
$row3[0,15] = 0.01
$row3[0,16] = 0.06
$row3[0,17] = 0.09
$row3[0,18] = 0.02
$row3[0,19] = 0.04
$row3[0,20] = 0
$row3[0,21] = 0
$row3[0,22] = 0
$row3[0,23] = 0
$row3[0,24] = 0
$row3[0,25] = 0
$row3[0,26] = 0.04
$row3[0,27] = 0.06
$row3[0,28] = 0.01
$row3[0,29] = 0.03
$row3[0,30] = 0
$row3[0,31] = 0
$row3[0,32] = 0
$row3[0,33] = 0
$row3[0,34] = 0
$row3[0,35] = 0
$row3[0,36] = 0.03
$row3[0,37] = 0.03
$row3[0,38] = 0
$row3[0,39] = 0.02
$row3[0,40] = 0
$row3[0,41] = 0
$row3[0,42] = 0
$row3[0,43] = 0
$row3[0,44] = 0
$row3[0,45] = 0
$row3[0,46] = 0.01
$row3[0,47] = 0.02
$row3[0,48] = 0
$row3[0,49] = 0.02
$ws.Range("A3:AX3").Value = $row3

$row4 = New-Object 'object[,]' 1,50
$row4[0,0] = 0.04
$row4[0,1] = 0.03
$row4[0,2] = 0.03
$row4[0,3] = 0.03
$row4[0,4] = 0.02
$row4[0,5] = 0.01
$row4[0,6] = 0.03
$row4[0,7] = 0.01
$row4[0,8] = 0
$row4[0,9] = 0
$row4[0,10] = 0.01
$row4[0,11] = 0.01
$row4[0,12] = 0.01
$row4[0,13] = 0.01
$row4[0,14] = 0.01
$row4[0,15] = 0
$row4[0,16] = 0.02
$row4[0,17] = 0.01
$row4[0,18] = 0
$row4[0,19] = 0
$row4[0,20] = 0.02
$row4[0,21] = 0.01
$row4[0,22] = 0.01
$row4[0,23] = 0.01
$row4[0,24] = 0.01
$row4[0,25] = 0
$row4[0,26] = 0
$row4[0,27] = 0
$row4[0,28] = 0
$row4[0,29] = 0
$row4[0,30] = 0.02
$row4[0,31] = 0.01
$row4[0,32] = 0.01
$row4[0,33] = 0.01
$row4[0,34] = 0.01
$row4[0,35] = 0.01
$row4[0,36] = 0
$row4[0,37] = 0
$row4[0,38] = 0
$row4[0,39] = 0
$row4[0,40] = 0.01
$row4[0,41] = 0.01
$row4[0,42] = 0.01
$row4[0,43] = 0.01
$row4[0,44] = 0
$row4[0,45] = 0
$row4[0,46] = 0
$row4[0,47] = 0
$row4[0,48] = 0
$row4[0,49] = 0
$ws.Range("A4:AX4").Value = $row4

$row5 = New-Object 'object[,]' 1,50
$row5[0,0] = 0.14
$row5[0,1] = 0.14
$row5[0,2] = 0.11
$row5[0,3] = 0.11
$row5[0,4] = 0.06
$row5[0,5] = 0.06
$row5[0,6] = 0.07
$row5[0,7] = 0.04
$row5[0,8] = 0
$row5[0,9] = 0
$row5[0,10] = 0.06
$row5[0,11] = 0.06
$row5[0,12] = 0.04
$row5[0,13] = 0.04
$row5[0,14] = 0.03
$row5[0,15] = 0.03
$row5[0,16] = 0.08
$row5[0,17] = 0.06
$row5[0,18] = 0
$row5[0,19] = 0
$row5[0,20] = 0.09
$row5[0,21] = 0.08
$row5[0,22] = 0.07
$row5[0,23] = 0.07
$row5[0,24] = 0.04
$row5[0,25] = 0.04
$row5[0,26] = 0.04
$row5[0,27] = 0.04
$row5[0,28] = 0
$row5[0,29] = 0
$row5[0,30] = 0.09
$row5[0,31] = 0.08
$row5[0,32] = 0.07
$row5[0,33] = 0.07
$row5[0,34] = 0.05
$row5[0,35] = 0.05
$row5[0,36] = 0
$row5[0,37] = 0
$row5[0,38] = 0
$row5[0,39] = 0
$row5[0,40] = 0.03
$row5[0,41] = 0.03
$row5[0,42] = 0.03
$row5[0,43] = 0.03
$row5[0,44] = 0
$row5[0,45] = 0
$row5[0,46] = 0
$row5[0,47] = 0
$row5[0,48] = 0
$row5[0,49] = 0
$ws.Range("A5:AX5").Value = $row5

